$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$row = 74
$ws.Cells.Item($row, 1).Value = 45803.93172453704
$ws.Cells.Item($row, 2).Value = 11
$ws.Cells.Item($row, 3).Value = 6
$ws.Cells.Item($row, 4).Value = 385
$ws.Cells.Item($row, 5).Value = 606
$ws.Cells.Item($row, 6).Value = 578
$ws.Cells.Item($row, 7).Value = 682
$ws.Cells.Item($row, 8).Value = 5472
$ws.Cells.Item($row, 9).Value = 682
$ws.Cells.Item($row, 10).Value = 2
$ws.Cells.Item($row, 11).Value = 2
$ws.Cells.Item($row, 12).Value = 673
$ws.Cells.Item($row, 13).Value = 30
$ws.Cells.Item($row, 14).Value = 5635
$ws.Cells.Item($row, 15).Value = 6862
